$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in / clear some previously-missing (imputed) values in column E
$ws.Range("E19").Value = -6.5
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7

# Remove the "RM 232" and "SC 92" rows entirely (data points dropped)
$ws.Rows("28:28").Delete()
$ws.Rows("26:26").Delete()

# After the row shift: SC 101 (now row 27) loses its E value,
# and SC 232 (now row 33) gains one.
$ws.Range("E27").ClearContents()
$ws.Range("E33").Value = -10.7
